# Refactor_8 remove code hardcode
# Populate the previously hard-coded single PatientID header sheet with
# the actual PEP patient IDs (rows 2 and 3), extending the used range
# from A1 to A1:A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "PEP_ID-2006019"
$ws.Range("A3").Value = "PEP_ID-2006022"
